{"js": "// Add a new \"Dasha.aboukalil@gmail.com\" entry after the last paragraph\n// (\"APPLE ACCOUNT: ... : Jamil1975\"): a blank spacer paragraph, then a\n// paragraph with a mailto hyperlink followed by \": kilmit sir: Dashjim1975\".\n\nconst body = context.document.body;\n\n// Find the last paragraph currently in the document body.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Blank paragraph acting as a spacer, matching the existing document style.\nconst spacerParagraph = lastParagraph.insertParagraph(\"\", \"After\");\n\n// New paragraph that will hold the hyperlink + trailing text.\nconst newParagraph = spacerParagraph.insertParagraph(\"Dasha.aboukalil@gmail.com\", \"After\");\n\n// Turn the whole paragraph's text into a mailto hyperlink (this also applies\n// the built-in \"Hyperlink\" character style, same as the existing hyperlink\n// earlier in the document).\nconst linkRange = newParagraph.getRange();\nlinkRange.hyperlink = \"mailto:Dasha.aboukalil@gmail.com\";\n\n// Append the trailing plain-text run after the hyperlink.\nnewParagraph.insertText(\": kilmit sir: Dashjim1975\", \"End\");\n\nawait context.sync();\n", "ps1": "# Add a new \"Dasha.aboukalil@gmail.com\" entry after the last paragraph\n# (\"APPLE ACCOUNT: ... : Jamil1975\"): a blank spacer paragraph, then a\n# paragraph with a mailto hyperlink followed by \": kilmit sir: Dashjim1975\".\n\n$d = $word.ActiveDocument\n\n# Locate the last paragraph currently in the document.\n$count = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($count)\n$r = $lastParagraph.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n# Blank spacer paragraph, matching the existing document style.\n$spacer = $d.Paragraphs.Last\n$spacer.Range.InsertParagraphAfter()\n\n# New paragraph that will hold the hyperlink + trailing text.\n$newParagraph = $d.Paragraphs.Last\n$newRange = $newParagraph.Range\n$newRange.Text = \"Dasha.aboukalil@gmail.com\"\n\n# Turn the paragraph's text into a mailto hyperlink (this also applies the\n# built-in \"Hyperlink\" character style, same as the existing hyperlink\n# earlier in the document).\n$newRange.Hyperlink = \"mailto:Dasha.aboukalil@gmail.com\"\n\n# Append the trailing plain-text run after the hyperlink.\n$tailRange = $d.Paragraphs.Last.Range\n$tailRange.Collapse(0)\n$tailRange.InsertAfter(\": kilmit sir: Dashjim1975\")\n"}
